# === LinkedIn workbook update: split Sheet1 into "Algorithm" + new "Design" sheet ===

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Algorithm"

# --- Row 62 on Algorithm: re-style like row 2 (blue "solved with note" look + hyperlink note) ---
$ws1.Range("B62:C62").Font.Color = 15773696
$ws1.Range("D62").Style = "Hyperlink"
$ws1.Range("D62").Font.Color = 15773696

# --- Create the new "Design" sheet right after "Algorithm" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Design"
$ws2.Columns.Item(1).ColumnWidth = 51.6

# --- Header + data rows (written in original authoring order) ---
$ws2.Range("A1").Value = "Description"
$ws2.Range("B1").Value = "Freq"
$ws2.Range("C1").Value = "Link"
$ws2.Range("A3").Value = "shorten url的各‍‌‌‍‍‍‌‌‌‌‌‌‍‍种变种，比如新添feature：click stats，就是统计每个short url被read多少次"
$ws2.Range("B3").Value = 9
$ws2.Range("A4").Value = "设计一个K/V store，支持基于单个KEY的insert/update/delete/fetch 操作，基本上照着RocksDB/Couchbase的实现来聊的"
$ws2.Range("B4").Value = 4
$ws2.Range("A5").Value = "设计一个metric系统，包括怎么收集，aggregation，存储，查询，dashboard，alert"
$ws2.Range("B5").Value = 3
$ws2.Range("A6").Value = "design monitoring system"
$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = "https://www.1point3acres.com/bbs/thread-542321-1-1.html"
$ws2.Range("A7").Value = "Design top k exception in 24 hours."
$ws2.Range("B7").Value = 14
$ws2.Range("A8").Value = "设计日历"
$ws2.Range("B8").Value = 8
$ws2.Range("D6").Value = "https://www.1point3acres.com/bbs/thread-532557-1-1.html"
$ws2.Range("C7").Value = "https://www.1point3acres.com/bbs/thread-531929-1-1.html"
$ws2.Range("C8").Value = "https://www.1point3acres.com/bbs/thread-531528-1-1.html "
$ws2.Range("A9").Value = "设计一个全球范围内的blacklist service"
$ws2.Range("B9").Value = 5
$ws2.Range("C9").Value = "https://www.1point3acres.com/bbs/thread-531259-1-1.html"
$ws2.Range("A10").Value = "设计distributed Logging System"
$ws2.Range("B10").Value = 2
$ws2.Range("A11").Value = "设计一个基于内存的streaming系统，stream以(timestamp‍‌‌‍‍‍‌‌‌‌‌‌‍‍, binary_size)的消息进入，然后client会query以ts结束大小为k的内容。"
$ws2.Range("B11").Value = 1
$ws2.Range("A13").Value = "存储在线用户的在网站上的活动"
$ws2.Range("B13").Value = 2
$ws2.Range("C13").Value = "https://www.1point3acres.com/bbs/thread-520850-1-1.html"
$ws2.Range("A14").Value = "Delayed Task Scheduler"
$ws2.Range("B14").Value = 5
$ws2.Range("A15").Value = "设计一个系统监督和管理领英第三方API的流量"
$ws2.Range("B15").Value = 1
$ws2.Range("A16").Value = "设计一个诊断系统，类似地理说的Kafka加上ag‍‌‌‍‍‍‌‌‌‌‌‌‍‍gregator的设计方式"
$ws2.Range("B16").Value = 2
$ws2.Range("A17").Value = "设计二级好友三级好友"
$ws2.Range("B17").Value = 3
$ws2.Range("A18").Value = "Amazon Product Page. 分析表之间的关系"
$ws2.Range("B18").Value = 2
$ws2.Range("C3").Value = "https://www.1point3acres.com/bbs/thread-495284-1-1.html"
$ws2.Range("A19").Value = "Design Hangman Game"
$ws2.Range("B19").Value = 2
$ws2.Range("C4").Value = "https://www.1point3acres.com/bbs/thread-492225-1-1.html "
$ws2.Range("A12").Value = "设计trending linkedin share post"
$ws2.Range("B12").Value = 2
$ws2.Range("D4").Value = "value体积比较大需要放在硬盘里面 另外随机写到硬盘会比较慢所以assume你要appending only"
$ws2.Range("A20").Value = "document repository"
$ws2.Range("B20").Value = 1
$ws2.Range("C14").Value = "https://soulmachine.gitbooks.io/system-design/cn/task-scheduler.html"
$ws2.Range("A2").Value = "https://www.1point3acres.com/bbs/thread-446923-1-1.html "

# --- Hyperlinks (same URL as the cell text) + re-apply the standard Hyperlink style ---
$url = "https://www.1point3acres.com/bbs/thread-542321-1-1.html"
$ws2.Hyperlinks.Add($ws2.Range("C6"), $url) | Out-Null
$ws2.Range("C6").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-532557-1-1.html"
$ws2.Hyperlinks.Add($ws2.Range("D6"), $url) | Out-Null
$ws2.Range("D6").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-531929-1-1.html"
$ws2.Hyperlinks.Add($ws2.Range("C7"), $url) | Out-Null
$ws2.Range("C7").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-531528-1-1.html "
$ws2.Hyperlinks.Add($ws2.Range("C8"), $url) | Out-Null
$ws2.Range("C8").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-531259-1-1.html"
$ws2.Hyperlinks.Add($ws2.Range("C9"), $url) | Out-Null
$ws2.Range("C9").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-520850-1-1.html"
$ws2.Hyperlinks.Add($ws2.Range("C13"), $url) | Out-Null
$ws2.Range("C13").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-495284-1-1.html"
$ws2.Hyperlinks.Add($ws2.Range("C3"), $url) | Out-Null
$ws2.Range("C3").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-492225-1-1.html "
$ws2.Hyperlinks.Add($ws2.Range("C4"), $url) | Out-Null
$ws2.Range("C4").Style = "Hyperlink"
$url = "https://soulmachine.gitbooks.io/system-design/cn/task-scheduler.html"
$ws2.Hyperlinks.Add($ws2.Range("C14"), $url) | Out-Null
$ws2.Range("C14").Style = "Hyperlink"
$url = "https://www.1point3acres.com/bbs/thread-446923-1-1.html "
$ws2.Hyperlinks.Add($ws2.Range("A2"), $url) | Out-Null
$ws2.Range("A2").Style = "Hyperlink"

# --- Selections / active sheet, matching the final saved view state ---
$ws1.Range("D54").Select()
$ws2.Range("A29").Select()
$ws2.Activate()
